# Apply latest Kimai data pull for pay period 21
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Dharam Pal (row 6): paid hours increased from 80 to 88
$ws.Range("F6").Value = 88
$ws.Range("G6").Value = 44
$ws.Range("H6").Value = 44
$ws.Range("I6").Value = "Extra 8.00 hours carry over"

# Yulia McCoy (row 10): paid hours increased from 50.25 to 54.25
$ws.Range("F10").Value = 54.25
$ws.Range("H10").Value = 44.25
$ws.Range("I10").Value = "Extra 4.25 hours carry over"
